$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.435.33"
$ws.Range("E2").Value = "  +5.27%  "
$ws.Range("D3").Value = "1.814.28"
$ws.Range("E3").Value = "  +5.06%  "
$ws.Range("E4").Value = "  -0.50%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9951"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.64%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5659"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +17.65%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3838"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +10.13%  "
$ws.Range("B9").Value = "OKB"
$ws.Range("C9").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "43.32"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07635"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.136"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +8.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.38"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +7.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.9938"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.87%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.239"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.14%  "
$ws.Range("D15").Value = "1.803.53"
$ws.Range("E15").Value = "  +4.41%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.232"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.55%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "92.32"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +6.25%  "
$ws.Range("E18").Value = "  +4.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06499"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.55%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9951"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.62%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.28"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.993"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.97%  "
$ws.Range("D23").Value = "28.456.80"
$ws.Range("E23").Value = "  +5.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.30"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.32%  "
$ws.Range("E25").Value = "  +1.42%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "20.81"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.94"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.88%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.378"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +14.77%  "
$ws.Range("D29").Value = "2.012.61"
$ws.Range("E29").Value = "  +4.72%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "123.70"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.35%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.150"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +10.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1047"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +12.46%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.771"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.00%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.621"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.76%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02316"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.09%  "
$ws.Range("B36").Value = "Algorand"
$ws.Range("C36").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2136"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.671"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +15.02%  "
$ws.Range("B38").Value = "Aptos"
$ws.Range("C38").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "11.64"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6429"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.45%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.049"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.95%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.06058"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.22%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9953"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.154"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.378"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.43%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.44"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.19%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5983"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.46%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.706"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.49%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "122.30"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.04%  "
$ws.Range("E49").Value = "  +4.93%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.142"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.24%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06841"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.06%  "
